$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.122.07'
$ws.Range('E2').Value = '  -0.73%  '
$ws.Range('D3').Value = '2.315.12'
$ws.Range('E3').Value = '  -1.86%  '
$ws.Range('E4').Value = '  +0.08%  '
$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '532.48'
$ws.Range('D5').Style = $origStyle
$ws.Range('E5').Value = '  +2.02%  '
$origStyle = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '132.19'
$ws.Range('D6').Style = $origStyle
$ws.Range('E6').Value = '  -3.19%  '
$ws.Range('E7').Value = '  -0.20%  '
$origStyle = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.536'
$ws.Range('D8').Style = $origStyle
$ws.Range('E8').Value = '  -0.52%  '
$ws.Range('D9').Value = '2.337.71'
$ws.Range('E9').Value = '  -1.66%  '
$origStyle = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.101'
$ws.Range('D10').Style = $origStyle
$ws.Range('E10').Value = '  -1.48%  '
$ws.Range('E11').Value = '  +0.21%  '
$ws.Range('E12').Value = '  -2.90%  '
$ws.Range('E13').Value = '  +0.56%  '
$ws.Range('D14').Value = '2.732.09'
$ws.Range('E14').Value = '  -1.70%  '
$origStyle = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '23.40'
$ws.Range('D15').Style = $origStyle
$ws.Range('E15').Value = '  -4.13%  '
$ws.Range('D16').Value = '57.184.60'
$ws.Range('E16').Value = '  -0.77%  '
$ws.Range('E17').Value = '  -2.75%  '
$ws.Range('D18').Value = '2.339.50'
$ws.Range('E18').Value = '  -0.93%  '
$origStyle = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '337.22'
$ws.Range('D19').Style = $origStyle
$ws.Range('E19').Value = '  +2.19%  '
$origStyle = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.42'
$ws.Range('D20').Style = $origStyle
$ws.Range('E20').Value = '  -2.09%  '
$origStyle = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.88'
$ws.Range('D21').Style = $origStyle
$ws.Range('E21').Value = '  +2.20%  '
$ws.Range('E22').Value = '  -2.34%  '
$ws.Range('E23').Value = '  -0.11%  '
$ws.Range('E24').Value = '  +0.18%  '
$ws.Range('E25').Value = '  +0.20%  '
$origStyle = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.68'
$ws.Range('D26').Style = $origStyle
$ws.Range('E26').Value = '  +4.28%  '
$origStyle = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.990'
$ws.Range('D27').Style = $origStyle
$ws.Range('E27').Value = '  -0.48%  '
$origStyle = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.32'
$ws.Range('D28').Style = $origStyle
$ws.Range('E28').Value = '  +0.01%  '
$origStyle = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '170.59'
$ws.Range('D29').Style = $origStyle
$ws.Range('E29').Value = '  +0.12%  '
$ws.Range('E30').Value = '  +0.76%  '
$ws.Range('D31').Value = '0.0₃0723'
$ws.Range('E31').Value = '  -3.28%  '
$origStyle = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.08'
$ws.Range('D32').Style = $origStyle
$ws.Range('E32').Value = '  -3.62%  '
$origStyle = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.52'
$ws.Range('D33').Style = $origStyle
$ws.Range('E33').Value = '  -0.45%  '
$origStyle = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.993'
$ws.Range('D35').Style = $origStyle
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('E36').Value = '  -4.06%  '
$origStyle = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.98'
$ws.Range('D37').Style = $origStyle
$ws.Range('E37').Value = '  -1.79%  '
$ws.Range('E38').Value = '  -3.15%  '
$ws.Range('E39').Value = '  +0.05%  '
$origStyle = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '39.17'
$ws.Range('D40').Style = $origStyle
$ws.Range('E40').Value = '  +1.30%  '
$origStyle = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '148.54'
$ws.Range('D41').Style = $origStyle
$ws.Range('E41').Value = '  -2.16%  '
$ws.Range('E42').Value = '  -1.72%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$origStyle = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.60'
$ws.Range('D43').Style = $origStyle
$ws.Range('E43').Value = '  -1.81%  '
$ws.Range('B44').Value = 'Bittensor'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$origStyle = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '280.48'
$ws.Range('D44').Style = $origStyle
$ws.Range('E44').Value = '  -1.32%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$origStyle = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.06'
$ws.Range('D45').Style = $origStyle
$ws.Range('E45').Value = '  -5.20%  '
$ws.Range('E46').Value = '  -1.38%  '
$ws.Range('E47').Value = '  -1.90%  '
$origStyle = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '18.67'
$ws.Range('D48').Style = $origStyle
$ws.Range('E48').Value = '  +1.37%  '
$origStyle = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.556'
$ws.Range('D49').Style = $origStyle
$ws.Range('E49').Value = '  -1.79%  '
$origStyle = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0216'
$ws.Range('D50').Style = $origStyle
$ws.Range('E50').Value = '  -2.32%  '
$ws.Range('E51').Value = '  -0.29%  '
